$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2025-12-12 Friday" "2025-12-13 Saturday"

Replace-Text "924×8=7392" "379×3=1137"
Replace-Text "157×3=471" "792×4=3168"
Replace-Text "136×6=816" "732×6=4392"
Replace-Text "708×7=4956" "900×5=4500"
Replace-Text "589×3=1767" "349×3=1047"

Replace-Text "722×3=2166" "585×4=2340"
Replace-Text "776×7=5432" "618×9=5562"
Replace-Text "767×2=1534" "961×3=2883"
Replace-Text "308×7=2156" "390×7=2730"
Replace-Text "372×2=744" "905×6=5430"

Replace-Text "161×8=1288" "406×8=3248"
Replace-Text "803×6=4818" "238×8=1904"
Replace-Text "950×4=3800" "116×8=928"
Replace-Text "119×5=595" "992×7=6944"
Replace-Text "484×5=2420" "292×5=1460"

Replace-Text "949×4=3796" "929×6=5574"
Replace-Text "211×4=844" "518×7=3626"
Replace-Text "551×8=4408" "295×8=2360"
Replace-Text "791×9=7119" "927×4=3708"
Replace-Text "318×5=1590" "554×8=4432"

Replace-Text "168×7=1176" "225×7=1575"
Replace-Text "722×2=1444" "755×6=4530"
Replace-Text "763×4=3052" "211×7=1477"
Replace-Text "779×9=7011" "817×6=4902"
Replace-Text "390×9=3510" "212×2=424"
